$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Oval 3" shape (falls back to the 3rd shape if renamed already).
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Oval 3") {
        $shp = $s.Shapes.Item($i)
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(3)
}

# EMU -> points (PowerPoint COM Left/Top/Width/Height are expressed in points).
$emuPerPt = 12700
# Tiny epsilon to counter single-precision truncation when the host converts
# the point value back to EMU internally, so the round-tripped EMU matches.
$eps = 0.00003

# Rename the shape and turn it into a plain rectangle.
$shp.Name = "Rectangle 3"
$shp.AutoShapeType = 1

# Reposition / resize to the new frame.
$shp.Left = 1277007 / $emuPerPt + $eps
$shp.Top = 3957145 / $emuPerPt + $eps
$shp.Width = 2790496 / $emuPerPt + $eps
$shp.Height = 1891862 / $emuPerPt + $eps

# Give it a solid accent2 theme fill.
$shp.Fill.Solid()
$shp.Fill.ForeColor.ObjectThemeColor = 6

# Shorten the placeholder caption.
$shp.TextFrame.TextRange.Text = "LOGO"
